$d = $word.ActiveDocument

# Locate the sentence that needs the Community Control address inserted.
$anchor = $d.Content.Duplicate
$found = $anchor.Find.Execute("Defendant shall provide written notice to the Office of Community Control at least 10 days prior to leaving Ohio.")

if ($found) {
    $sentenceStart = $anchor.Start
    $sentenceEnd = $anchor.End

    # Work out exactly where "...Community Control" ends within the
    # sentence so the run can be split right after it (i.e. before
    # " at least 10 days...").
    $marker = $d.Range($sentenceStart, $sentenceEnd)
    $markerFound = $marker.Find.Execute("Community Control")
    $splitPoint = $marker.End

    # Pull out (and blank) the trailing part of the original sentence so it
    # can be re-inserted afterwards as its own run.
    $tail = $d.Range($splitPoint, $sentenceEnd)
    $tailText = $tail.Text
    $tail.Text = ""

    # Insert the new address clause as its own run right after
    # "Community Control".
    $addressPoint = $d.Range($splitPoint, $splitPoint)
    $addressPoint.InsertAfter(", located at 70 N. Union St., Delaware, OH 43015,")

    # Re-insert the original tail text (" at least 10 days prior to leaving
    # Ohio.") as its own run after the new address run.
    $tailPoint = $d.Range($addressPoint.End, $addressPoint.End)
    $tailPoint.InsertAfter($tailText)

    # Touch (no-op toggle) the formatting of each new piece so the engine
    # keeps them as distinct runs instead of silently re-merging them with
    # their neighbours -- the visible formatting ends up unchanged (same
    # Palatino Linotype / bCs / sz 20 / szCs 20 as the surrounding text).
    $addressPoint.Font.Bold = 1
    $addressPoint.Font.Bold = 0
    $tailPoint.Font.Bold = 1
    $tailPoint.Font.Bold = 0
}
